$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '24.767.70'
$ws.Range('E2').Value = '  +0.18%  '
# Row 3
$ws.Range('D3').Value = '1.705.37'
$ws.Range('E3').Value = '  +0.32%  '
# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9969'
$ws.Range('E4').Value = '  -0.68%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '318.34'
$ws.Range('E5').Value = '  +0.61%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9969'
$ws.Range('E6').Value = '  -0.56%  '
# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3928'
$ws.Range('E7').Value = '  -0.02%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4069'
$ws.Range('E8').Value = '  +0.37%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.497'
$ws.Range('E9').Value = '  -1.44%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '54.43'
$ws.Range('E10').Value = '  +2.96%  '
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9971'
$ws.Range('E11').Value = '  -0.77%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08823'
$ws.Range('E12').Value = '  -0.61%  '
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '26.38'
$ws.Range('E13').Value = '  +11.26%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.489'
$ws.Range('E14').Value = '  +0.77%  '
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.153'
$ws.Range('E15').Value = '  +0.34%  '
# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001357'
$ws.Range('E16').Value = '  +2.78%  '
# Row 17
$ws.Range('D17').Value = '1.693.46'
$ws.Range('E17').Value = '  -0.56%  '
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '97.79'
$ws.Range('E18').Value = '  -1.69%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07159'
$ws.Range('E19').Value = '  +1.46%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '20.57'
$ws.Range('E20').Value = '  +3.85%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.312'
$ws.Range('E21').Value = '  +3.47%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9949'
$ws.Range('E22').Value = '  -0.93%  '
# Row 23
$ws.Range('E23').Value = '  -2.28%  '
# Row 24
$ws.Range('D24').Value = '24.752.60'
$ws.Range('E24').Value = '  +0.15%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.024'
$ws.Range('E25').Value = '  -3.85%  '
# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.332'
$ws.Range('E26').Value = '  -0.79%  '
# Row 27
$ws.Range('E27').Value = '  +1.60%  '
# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '167.15'
$ws.Range('E28').Value = '  +1.68%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.969'
$ws.Range('E29').Value = '  +15.76%  '
# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.503'
$ws.Range('E30').Value = '  -3.87%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '145.10'
$ws.Range('E31').Value = '  +6.82%  '
# Row 32
$ws.Range('D32').Value = '1.883.79'
$ws.Range('E32').Value = '  -0.47%  '
# Row 33
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.08816'
$ws.Range('E33').Value = '  -2.33%  '
# Row 34
$ws.Range('B34').Value = 'WEMIXTOKEN'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.173'
$ws.Range('E34').Value = '  +10.90%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.073'
$ws.Range('E35').Value = '  +0.43%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.251'
$ws.Range('E36').Value = '  -5.16%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.03120'
$ws.Range('E37').Value = '  +3.48%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2812'
$ws.Range('E38').Value = '  +1.89%  '
# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.8457'
$ws.Range('E39').Value = '  +9.84%  '
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '10.91'
$ws.Range('E40').Value = '  -0.69%  '
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.09227'
$ws.Range('E41').Value = '  +0.05%  '
# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '14.22'
$ws.Range('E42').Value = '  -1.58%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.484'
$ws.Range('E43').Value = '  +0.89%  '
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '17.62'
$ws.Range('E44').Value = '  +9.22%  '
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.717'
$ws.Range('E45').Value = '  +5.08%  '
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.7468'
$ws.Range('E46').Value = '  +4.03%  '
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.277'
$ws.Range('E47').Value = '  +1.30%  '
# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.394'
$ws.Range('E48').Value = '  +2.67%  '
# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.9963'
$ws.Range('E49').Value = '  -0.67%  '
# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '140.78'
$ws.Range('E50').Value = '  +0.66%  '
# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.08250'
$ws.Range('E51').Value = '  +3.52%  '
